{"js": "// The author corrected the date on the title page from \"13 February 2019\"\n// to \"25 February 2019\". Re-typing just the day number is what produced the\n// run split (\"25\" | \" February 2019\") and moved Word's automatic \"_GoBack\"\n// (last-edit) bookmark from its previous location (right after\n// \"Slide & Connect\") to sit between the new day number and the rest of the\n// date.\n\n// 1) Drop the \"_GoBack\" bookmark from wherever it currently lives.\nconst oldGoBack = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\noldGoBack.load(\"isNullObject\");\nawait context.sync();\nif (!oldGoBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 2) Locate the date paragraph and split it so the leading \"13\" can be\n//    replaced independently of \" February 2019\".\nconst dateMatches = context.document.body.search(\"13 February 2019\", {\n  matchCase: true,\n  matchWholeWord: false\n});\ndateMatches.load(\"items\");\nawait context.sync();\n\nconst dateRange = dateMatches.items[0];\nconst dateParts = dateRange.split([\" \"], true, true, true);\ndateParts.load(\"items\");\nawait context.sync();\n\nconst dayRange = dateParts.items[0]; // \"13\"\n\n// 3) Replace \"13\" with \"25\".\ndayRange.insertText(\"25\", \"Replace\");\nawait context.sync();\n\n// 4) Re-create \"_GoBack\" right after the new day number (before\n//    \" February 2019\"), matching where Word leaves it after an in-place\n//    edit.\nconst afterDay = dayRange.getRange(\"End\");\nafterDay.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The author corrected the date on the title page from \"13 February 2019\"\n# to \"25 February 2019\". Re-typing just the day number is what produced the\n# run split (\"25\" | \" February 2019\") and moved Word's automatic \"_GoBack\"\n# (last-edit) bookmark from its previous location (right after\n# \"Slide & Connect\") to sit between the new day number and the rest of the\n# date.\n\n$d = $word.ActiveDocument\n\n# 1) Find the date line and narrow the range down to just the leading \"13\".\n$dateRange = $d.Content\n$dateRange.Find.Execute(\"13 February 2019\") | Out-Null\n$dayRange = $d.Range($dateRange.Start, $dateRange.Start + 2)\n\n# 2) Replace \"13\" with \"25\".\n$dayRange.Text = \"25\"\n\n# 3) Re-create \"_GoBack\" right after the new day number (before\n#    \" February 2019\"). Bookmark names are unique, so adding it again\n#    simply relocates the existing one away from \"Slide & Connect\".\n$insertionPoint = $d.Range($dayRange.End, $dayRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $insertionPoint) | Out-Null\n"}
